# Adds data for 2022-08-04: the rolling "Through YYYY-MM-DD" carjacking
# report moves forward one day (July 26 -> July 27) and a batch of
# carjacking counts across many months/neighborhoods are incremented or
# newly populated to reflect the latest pull from the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet and update the "through" date shown in the header cell.
$ws.Name = "Through 2022-07-27"
$ws.Range("B1").Value = "July 2022 (through July 27)"

# Updated / newly-populated monthly counts per neighborhood row.
$ws.Range("AD2").Value = 13
$ws.Range("P3").Value = 9
$ws.Range("W3").Value = 2
$ws.Range("I4").Value = 1
$ws.Range("AD5").Value = 5
$ws.Range("AK5").Value = 5
$ws.Range("AY5").Value = 3
$ws.Range("I6").Value = 4
$ws.Range("P8").Value = 20
$ws.Range("I12").Value = 3
$ws.Range("B13").Value = 2
$ws.Range("P16").Value = 3
$ws.Range("I17").Value = 1
$ws.Range("AD18").Value = 1
$ws.Range("B33").Value = 3
$ws.Range("AK33").Value = 1
$ws.Range("P38").Value = 2
$ws.Range("W38").Value = 1
$ws.Range("B39").Value = 3
$ws.Range("P44").Value = 1
$ws.Range("P47").Value = 1
$ws.Range("W47").Value = 2
$ws.Range("AY47").Value = 1
$ws.Range("AY49").Value = 1
$ws.Range("AK51").Value = 1
$ws.Range("AK52").Value = 4
$ws.Range("I56").Value = 1
$ws.Range("AR62").Value = 1
$ws.Range("B66").Value = 1
$ws.Range("I79").Value = 3
$ws.Range("P94").Value = 3
$ws.Range("I96").Value = 9
